$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("IRS-Cleared")
$ws2 = $wb.Worksheets.Item("IRS-Bilateral")
$ws1.Range("B2").Value = "ACUOSG8745"
